# Edit script: applies the two changes captured in the target diff.
#  1) Slide 6's table gets a new (built-in) table style applied.
#  2) The deck's theme (theme1.xml, used by the Slide Master) is recolored
#     from the "Integral" palette to the "Office Theme" palette - i.e. the
#     Design/Theme colour scheme that PowerPoint applies when you switch the
#     presentation's theme away from "Integral" to the default "Office Theme".

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 6 ------------------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{AB0D9AFE-BD5F-4B72-A724-0FE232F50CE6}")
    }
}

# --- 2) Re-theme the presentation's Slide Master colour scheme ---------------
# Colour order exposed by ThemeColorScheme.Item(n):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
#   9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGB values use the COLORREF (BGR) byte order PowerPoint's COM automation
# expects, so each target "Office Theme" RGB hex colour is byte-swapped.
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # dk1      000000
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = 0x6A5444   # dk2      44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$tcs.Item(9).RGB  = 0xC47244   # accent5  4472C4
$tcs.Item(10).RGB = 0x47AD70   # accent6  70AD47
$tcs.Item(11).RGB = 0xC16305   # hlink    0563C1
$tcs.Item(12).RGB = 0x724F95   # folHlink 954F72
